$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.142.13"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "3.122.83"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.122.73"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("D13").Value = "3.657.10"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("E14").Value = "  +3.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000165"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "58.206.74"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "3.122.75"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "342.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.514"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "0.0₃0924"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.06%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0669"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D44").Value = "3.159.23"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0263"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("D48").Value = "2.277.82"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  +4.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.34%  "
